$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (32 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3832.7917
$ws.Range("I113").Value = 3207.9167
$ws.Range("J113").Value = 4457.6665
$ws.Range("K113").Value = 3207.9167
$ws.Range("L113").Value = 4457.6665
$ws.Range("M113").Value = 46.08329999999978
$ws.Range("N113").Value = -10965.6665
$ws.Range("H116").Value = 2691.8235
$ws.Range("I116").Value = 2217.2222
$ws.Range("J116").Value = 3225.75
$ws.Range("K116").Value = 2217.2222
$ws.Range("L116").Value = 3225.75
$ws.Range("M116").Value = 1224.7778
$ws.Range("N116").Value = -10109.75
$ws.Range("H132").Value = 4278.7046
$ws.Range("I132").Value = 2033.079
$ws.Range("K132").Value = 6099.237
$ws.Range("M132").Value = -3569.237
$ws.Range("H137").Value = 1069273.2
$ws.Range("I137").Value = 786.0227
$ws.Range("J137").Value = 6945953
$ws.Range("K137").Value = 2358.0681
$ws.Range("L137").Value = 20837859
$ws.Range("M137").Value = 191.9319
$ws.Range("N137").Value = -20842959
$ws.Range("H138").Value = 2441430.8
$ws.Range("I138").Value = 1027
$ws.Range("J138").Value = 4003289
$ws.Range("K138").Value = 3081
$ws.Range("L138").Value = 12009867
$ws.Range("M138").Value = 2059
$ws.Range("N138").Value = -12020147

# --- Sheet: ARM (32 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4257.355
$ws.Range("I32").Value = 2764.8354
$ws.Range("K32").Value = 2764.8354
$ws.Range("M32").Value = -2477.8354
$ws.Range("H61").Value = 1463.02
$ws.Range("I61").Value = 1376.7188
$ws.Range("J61").Value = 1616.4445
$ws.Range("K61").Value = 1376.7188
$ws.Range("L61").Value = 1616.4445
$ws.Range("M61").Value = -1164.7188
$ws.Range("N61").Value = -2040.4445
$ws.Range("H122").Value = 2074.2354
$ws.Range("I122").Value = 1723.7778
$ws.Range("J122").Value = 2468.5
$ws.Range("K122").Value = 5171.3334
$ws.Range("L122").Value = 7405.5
$ws.Range("M122").Value = -2721.3334
$ws.Range("N122").Value = -12305.5
$ws.Range("H132").Value = 1679.9615
$ws.Range("I132").Value = 1380.0555
$ws.Range("J132").Value = 2354.75
$ws.Range("K132").Value = 4140.166499999999
$ws.Range("L132").Value = 7064.25
$ws.Range("M132").Value = -1610.166499999999
$ws.Range("N132").Value = -12124.25
$ws.Range("H136").Value = 1463.02
$ws.Range("I136").Value = 1376.7188
$ws.Range("J136").Value = 1616.4445
$ws.Range("K136").Value = 4130.1564
$ws.Range("L136").Value = 4849.333500000001
$ws.Range("M136").Value = -1580.1564
$ws.Range("N136").Value = -9949.333500000001

# --- Sheet: BSM (18 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1785.3846
$ws.Range("I99").Value = 1523.3334
$ws.Range("J99").Value = 2375
$ws.Range("K99").Value = 1523.3334
$ws.Range("L99").Value = 2375
$ws.Range("M99").Value = -25.33339999999998
$ws.Range("N99").Value = -5371
$ws.Range("H105").Value = 2322.5
$ws.Range("I105").Value = 2223.3333
$ws.Range("K105").Value = 2223.3333
$ws.Range("M105").Value = -476.3332999999998
$ws.Range("H134").Value = 648829.7
$ws.Range("I134").Value = 1084169.8
$ws.Range("J134").Value = 4526.44
$ws.Range("K134").Value = 3252509.4
$ws.Range("L134").Value = 13579.32
$ws.Range("M134").Value = -3249974.4
$ws.Range("N134").Value = -18649.32

# --- Sheet: CRP (39 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8066458.5
$ws.Range("I31").Value = 1241.4546
$ws.Range("J31").Value = 12502328
$ws.Range("K31").Value = 1241.4546
$ws.Range("L31").Value = 12502328
$ws.Range("M31").Value = -946.4546
$ws.Range("N31").Value = -12502918
$ws.Range("H34").Value = 8066458.5
$ws.Range("I34").Value = 1241.4546
$ws.Range("J34").Value = 12502328
$ws.Range("K34").Value = 1241.4546
$ws.Range("L34").Value = 12502328
$ws.Range("M34").Value = -1039.4546
$ws.Range("N34").Value = -12502732
$ws.Range("H99").Value = 2609245.5
$ws.Range("I99").Value = 3338199.2
$ws.Range("J99").Value = 5838.4287
$ws.Range("K99").Value = 3338199.2
$ws.Range("L99").Value = 5838.4287
$ws.Range("M99").Value = -3336701.2
$ws.Range("N99").Value = -8834.4287
$ws.Range("H126").Value = 2609245.5
$ws.Range("I126").Value = 3338199.2
$ws.Range("J126").Value = 5838.4287
$ws.Range("K126").Value = 10014597.6
$ws.Range("L126").Value = 17515.2861
$ws.Range("M126").Value = -10012127.6
$ws.Range("N126").Value = -22455.2861
$ws.Range("H132").Value = 1060404.1
$ws.Range("I132").Value = 2229.4814
$ws.Range("J132").Value = 4631743.5
$ws.Range("K132").Value = 6688.4442
$ws.Range("L132").Value = 13895230.5
$ws.Range("M132").Value = -4158.4442
$ws.Range("N132").Value = -13900290.5
$ws.Range("H134").Value = 2431.2163
$ws.Range("I134").Value = 2634.64
$ws.Range("K134").Value = 7903.92
$ws.Range("M134").Value = -5368.92

# --- Sheet: CUL (11 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 918.85
$ws.Range("J131").Value = 920.9693600000001
$ws.Range("L131").Value = 2762.90808
$ws.Range("N131").Value = -12842.90808
$ws.Range("H137").Value = 37596424
$ws.Range("I137").Value = 1932.5
$ws.Range("J137").Value = 67672020
$ws.Range("K137").Value = 5797.5
$ws.Range("L137").Value = 203016060
$ws.Range("M137").Value = -697.5
$ws.Range("N137").Value = -203026260

# --- Sheet: GSM (40 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4621.778
$ws.Range("I70").Value = 4425
$ws.Range("K70").Value = 4425
$ws.Range("M70").Value = -4155
$ws.Range("H73").Value = 4621.778
$ws.Range("I73").Value = 4425
$ws.Range("K73").Value = 4425
$ws.Range("M73").Value = -3489
$ws.Range("H80").Value = 2153.5
$ws.Range("I80").Value = 2113.125
$ws.Range("J80").Value = 2207.3333
$ws.Range("K80").Value = 2113.125
$ws.Range("L80").Value = 2207.3333
$ws.Range("M80").Value = -1115.125
$ws.Range("N80").Value = -4203.3333
$ws.Range("H83").Value = 2153.5
$ws.Range("I83").Value = 2113.125
$ws.Range("J83").Value = 2207.3333
$ws.Range("K83").Value = 10565.625
$ws.Range("L83").Value = 11036.6665
$ws.Range("M83").Value = -5573.625
$ws.Range("N83").Value = -21020.6665
$ws.Range("H93").Value = 18083.666
$ws.Range("J93").Value = 18083.666
$ws.Range("L93").Value = 18083.666
$ws.Range("N93").Value = -21827.666
$ws.Range("H122").Value = 80640.69500000001
$ws.Range("I122").Value = 93708.55
$ws.Range("J122").Value = 8767.5
$ws.Range("K122").Value = 281125.65
$ws.Range("L122").Value = 26302.5
$ws.Range("M122").Value = -278675.65
$ws.Range("N122").Value = -31202.5
$ws.Range("H132").Value = 2442090.5
$ws.Range("I132").Value = 2719.8
$ws.Range("J132").Value = 6253607
$ws.Range("K132").Value = 8159.400000000001
$ws.Range("L132").Value = 18760821
$ws.Range("M132").Value = -5629.400000000001
$ws.Range("N132").Value = -18765881

# --- Sheet: WVR (11 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 47279.69
$ws.Range("I122").Value = 10998.105
$ws.Range("K122").Value = 32994.315
$ws.Range("M122").Value = -30544.315
$ws.Range("H132").Value = 1876.3135
$ws.Range("I132").Value = 2118.7273
$ws.Range("J132").Value = 1412.5652
$ws.Range("K132").Value = 6356.1819
$ws.Range("L132").Value = 4237.6956
$ws.Range("M132").Value = -3826.1819
$ws.Range("N132").Value = -9297.695599999999
